$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "«FECHA_CARTA_ART_77__P_T_COMPR_ADJUD»" "27 de noviembre de 2019"
Replace-Text "«PROVEEDOR_1»" "COMERCIALIZADORA DE EQUIPOS Y SERVICIOS PARA LA INDUSTRIA AG, S.A. DE C.V."
Replace-Text "«REPRESENTANTE_LEGA_PROVEEDOR_1»" "ARMANDO ALTAMIRANO ALVARADO"
Replace-Text "«NOMBRE_PROCEDIMIENTO_»" "INV 78"

# The "«TIPO_DE_PROCEDIMIENTO_INV_CON_DIR_y_NO»" merge field result is split
# across two runs with a proofErr spell-check tag between them ("«" and
# "TIPO_DE_PROCEDIMIENTO_INV_CON_DIR_y_NO»"), and is immediately followed by
# "número" / proofErr / " " as three more runs. Scope the Find/Replace calls
# to the specific paragraph so the generic "número " text elsewhere in the
# document (already a single run) is left untouched.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "TIPO_DE_PROCEDIMIENTO_INV_CON_DIR_y_NO") {
        $p.Range.Find.Execute("«TIPO_DE_PROCEDIMIENTO_INV_CON_DIR_y_NO»", $true, $false, $false, $false, $false, $true, 1, $false, "INVITACION A CUANDO MENOS TRES PERSONAS", 2)
        $p.Range.Find.Execute("número ", $true, $false, $false, $false, $false, $true, 1, $false, "número ", 2)
        break
    }
}

$sec = $d.Sections(1)
$sec.PageSetup.TopMargin = 2552 / 20.0
$sec.PageSetup.RightMargin = 1418 / 20.0
$sec.PageSetup.BottomMargin = 1418 / 20.0
$sec.PageSetup.LeftMargin = 1418 / 20.0
$sec.PageSetup.HeaderDistance = 425 / 20.0
$sec.PageSetup.FooterDistance = 442 / 20.0
